# Updated notebook, reran simulation
# - Two new HKL reflections "Holden" and "Rizzie Spiral" were added to the
#   averaging table (inserted right after "Spiral5"), each with freshly
#   simulated ratio data across the 18 parameter columns.
# - "Thomas Hex" was renamed to "Matthies Hex".
# - Because two rows were inserted, the table grew from 29 to 31 data rows
#   (dimension A1:T29 -> A1:T31); everything that used to follow "Spiral5"
#   shifted down by two rows, carrying its original data with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the old row 4 (i.e. right after the
# "Spiral5" row), shifting all the existing "RotRing OmegaMax-90"..."Michael-
# SNHex" rows down by two (rows 4-29 -> rows 6-31) along with their data.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Give the new A4/A5 index cells the same bold/centered/bordered look as the
# rest of column A (style index 1 in the original file).
$ws.Range("A4:A5").Font.Bold = $true
$ws.Range("A4:A5").HorizontalAlignment = -4108
$ws.Range("A4:A5").VerticalAlignment = -4160
$ws.Range("A4:A5").Borders.LineStyle = 1

# Row 4: new "Holden" reflection
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.038241207505334
$ws.Range("D4").Value = 0.9389187403456561
$ws.Range("E4").Value = 0.9389187403456561
$ws.Range("F4").Value = 1.003298826768173
$ws.Range("G4").Value = 1.013486107553851
$ws.Range("H4").Value = 1.006255349760796
$ws.Range("I4").Value = 0.9722054178907777
$ws.Range("J4").Value = 1.038241207505334
$ws.Range("K4").Value = 1.038241207505334
$ws.Range("L4").Value = 1.003298826768173
$ws.Range("M4").Value = 0.9711087835569143
$ws.Range("N4").Value = 0.9711087835569143
$ws.Range("O4").Value = 0.9714743283348688
$ws.Range("P4").Value = 0.9934862582063877
$ws.Range("Q4").Value = 0.9934862582063877
$ws.Range("R4").Value = 1.004674995531124
$ws.Range("S4").Value = 1.004674995531124
$ws.Range("T4").Value = 0.9954009416374313

# Row 5: new "Rizzie Spiral" reflection
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 0.9075177189352477
$ws.Range("D5").Value = 1.560709027560744
$ws.Range("E5").Value = 1.560709027560744
$ws.Range("F5").Value = 0.8379125430050716
$ws.Range("G5").Value = 0.8582055962772724
$ws.Range("H5").Value = 0.5817437372523757
$ws.Range("I5").Value = 1.31044248474102
$ws.Range("J5").Value = 0.9075177189352477
$ws.Range("K5").Value = 0.9075177189352477
$ws.Range("L5").Value = 0.8379125430050716
$ws.Range("M5").Value = 1.199310785282908
$ws.Range("N5").Value = 1.199310785282908
$ws.Range("O5").Value = 1.236354685102279
$ws.Range("P5").Value = 1.102046429833688
$ws.Range("Q5").Value = 1.102046429833688
$ws.Range("R5").Value = 1.053414252109078
$ws.Range("S5").Value = 1.053414252109078
$ws.Range("T5").Value = 1.009421851295289

# Rename "Thomas Hex" -> "Matthies Hex" (now on row 11 after the shift).
$ws.Range("B11").Value = "Matthies Hex"
